# The workbook's single sheet ("Artfynd") holds one species-observation per
# data row (rows 2-13, columns A:AY). The edit described by the diff is a
# cyclic rotation of those 12 data rows: the record that used to sit in row 2
# moves down to become the last record (row 13), and every other record
# (originally rows 3-13) shifts up by one row (row 3 -> row 2, row 4 -> row 3,
# ..., row 13 -> row 12). The header row (row 1) is untouched.
#
# We reproduce that with whole-row Range.Copy operations, stashing the
# soon-to-be-displaced row 2 in a scratch row far below the used range, then
# shifting rows 3..13 up one at a time, and finally dropping the stashed row
# into row 13. Each destination is explicitly cleared immediately before the
# copy so that columns which are populated in the source but blank in the
# destination (or vice versa) end up correct rather than retaining stale
# leftovers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 13
$lastCol = "AY"
$scratchRow = 1000

$srcRange = "A" + $firstDataRow + ":" + $lastCol + $firstDataRow
$scratchRange = "A" + $scratchRow + ":" + $lastCol + $scratchRow

# 1) Stash row 2 (it will become the new row 13).
$ws.Range($scratchRange).ClearContents()
$ws.Range($srcRange).Copy($ws.Range($scratchRange))

# 2) Shift rows 3..13 up into rows 2..12, top to bottom.
for ($r = $firstDataRow + 1; $r -le $lastDataRow; $r++) {
    $dstRow = $r - 1
    $dstRange = "A" + $dstRow + ":" + $lastCol + $dstRow
    $curRange = "A" + $r + ":" + $lastCol + $r

    $ws.Range($dstRange).ClearContents()
    $ws.Range($curRange).Copy($ws.Range($dstRange))
}

# 3) Drop the stashed original row 2 into row 13.
$finalRange = "A" + $lastDataRow + ":" + $lastCol + $lastDataRow
$ws.Range($finalRange).ClearContents()
$ws.Range($scratchRange).Copy($ws.Range($finalRange))

# 4) Clean up the scratch row so it doesn't extend the used range.
$ws.Range($scratchRange).ClearContents()
